$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: merge a paragraph's runs that start at $fromOffset (paragraph-local
# char offset) through the end of the paragraph text (excluding the
# paragraph mark) into the run that immediately precedes $fromOffset. This
# is done by capturing the text, clearing it (deleting it entirely) and
# then using InsertAfter on a collapsed range sitting right at the
# boundary -- which extends the preceding run rather than minting a
# brand-new one, so unrelated run-level attributes (e.g. w:rsidRPr) on
# unchanged text survive.
# ---------------------------------------------------------------------------
function Merge-Tail($paraIndex, $fromOffset) {
    $p = $d.Paragraphs($paraIndex).Range
    $splitPos = $p.Start + $fromOffset
    $tail = $d.Range($splitPos, $p.End - 1)
    $txt = $tail.Text
    $tail.Text = ""
    $anchor = $d.Range($splitPos, $splitPos)
    $anchor.InsertAfter($txt)
}

# ---------------------------------------------------------------------------
# Helper: force a run split at an absolute document character offset by
# dropping a temporary bookmark there and immediately deleting it again.
# The bookmark insertion/removal leaves no trace, but the run boundary it
# creates persists.
# ---------------------------------------------------------------------------
function Split-At($absPos) {
    $bm = $d.Bookmarks.Add("__tmp_split__", $d.Range($absPos, $absPos))
    $d.Bookmarks("__tmp_split__").Delete()
}

# === Paragraph 1: "3.1 Folder Lock (Integrity)" -> "...(Availability)" =====
# split into 3 runs: "3.1 Folder Lock (" / "Availability" / ")"
$p1 = $d.Paragraphs(1).Range
$prefixLen1 = "3.1 Folder Lock (".Length
$splitPos1 = $p1.Start + $prefixLen1
$tail1 = $d.Range($splitPos1, $p1.End - 1)
$tail1.Text = ""
$anchor1 = $d.Range($splitPos1, $splitPos1)
$anchor1.InsertAfter("Availability)")
$availLen = "Availability".Length
$splitPos1b = $splitPos1 + $availLen
Split-At $splitPos1
Split-At $splitPos1b

# === Paragraph 2: merge "3.2 DES Encrypt" + "ion Procedure (...)" =========
$off2 = "3.2 DES Encrypt".Length
Merge-Tail 2 $off2

# === Paragraph 3: "3.3 File Hashing (Availability)" -> "...(Integrity)" ===
# split into 3 runs: "3.3 File Hashing (" / "Integrity" / ")"
$p3 = $d.Paragraphs(3).Range
$prefixLen3 = "3.3 File Hashing (".Length
$splitPos3 = $p3.Start + $prefixLen3
$tail3 = $d.Range($splitPos3, $p3.End - 1)
$tail3.Text = ""
$anchor3 = $d.Range($splitPos3, $splitPos3)
$anchor3.InsertAfter("Integrity)")
$integrityLen = "Integrity".Length
$splitPos3b = $splitPos3 + $integrityLen
Split-At $splitPos3
Split-At $splitPos3b

# === Paragraph 4: merge "3.4 Password Stren" + "gth Procedure (...)" ======
$off4 = "3.4 Password Stren".Length
Merge-Tail 4 $off4

# === Paragraph 5: merge "3.5 Email Veri" + "fication Code (...)" ==========
$off5 = "3.5 Email Veri".Length
Merge-Tail 5 $off5

# === Paragraph 6: resplit "3.6 Captcha (Confidentiality)" & move bookmark =
# original split after "3.6 Captch" (10 chars); target split after
# "3.6 Captcha (Co" (15 chars), with the _GoBack bookmark sitting there.
$off6 = "3.6 Captch".Length
Merge-Tail 6 $off6
$p6 = $d.Paragraphs(6).Range
$newOff6 = "3.6 Captcha (Co".Length
$newSplit6 = $p6.Start + $newOff6
Split-At $newSplit6
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($newSplit6, $newSplit6))

Write-Output "done"
